# Updated vignette tables and figure caption format
#
# Renames the per-sample column headers in column A (rows 2-9) of the
# "asv_small" demo table from the old V1_x_t0 / V2_x_t0 naming scheme to the
# simpler Sample1..Sample8 scheme, and moves the sheet's active-cell
# selection to C17 (matching where the author's cursor ended up before
# saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sample1"
$ws.Range("A3").Value = "Sample2"
$ws.Range("A4").Value = "Sample3"
$ws.Range("A5").Value = "Sample4"
$ws.Range("A6").Value = "Sample5"
$ws.Range("A7").Value = "Sample6"
$ws.Range("A8").Value = "Sample7"
$ws.Range("A9").Value = "Sample8"

$ws.Range("C17").Select() | Out-Null
